# Weekly update: insert a new price record as row 80 for
# "Fruta, Feria Lagunitas de Puerto Montt - Uva", pushing the existing
# rows 80..149 down to 81..150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 80 (shifts 80-149 -> 81-150)
$ws.Range("A80").EntireRow.Insert()

# Populate the newly inserted row 80 with the new data record
$ws.Range("A80").Value = 4
$ws.Range("B80").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C80").Value = "Los Lagos"
$ws.Range("D80").Value = 44512
$ws.Range("E80").Value = 10
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100109
$ws.Range("H80").Value = "Uva"
$ws.Range("I80").Value = 100109001
$ws.Range("J80").Value = "Uva"
$ws.Range("K80").Value = "Superior Seedless"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 200
$ws.Range("N80").Value = 32000
$ws.Range("O80").Value = 33000
$ws.Range("P80").Value = 32500
$ws.Range("Q80").Value = "`$/bandeja 10 kilos"
$ws.Range("R80").Value = "Provincia de Limar" + [char]0x00ED
$ws.Range("S80").Value = 3250
$ws.Range("T80").Value = 10
